# Updated cryptos list on Thu Mar 16 20:31:52 UTC 2023 with GitHub Actions
#
# Price (D) and Volume(1h) (E) columns are stored as plain text in the
# workbook (e.g. "24.984.52", "  +2.32%  "), so values are entered with a
# leading apostrophe to force text storage and avoid Excel's automatic
# number/date coercion (which would also collapse formatting such as
# trailing zeros or thousands separators).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.984.52"
$ws.Range("E2").Value = "'  +2.32%  "
$ws.Range("D3").Value = "'1.680.02"
$ws.Range("E3").Value = "'  +1.83%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "'  -0.15%  "
$ws.Range("D5").Value = "'328.86"
$ws.Range("E5").Value = "'  +6.99%  "
$ws.Range("D6").Value = "'0.9995"
$ws.Range("D7").Value = "'0.3665"
$ws.Range("E7").Value = "'  +1.30%  "
$ws.Range("D8").Value = "'47.04"
$ws.Range("E8").Value = "'  -0.95%  "
$ws.Range("D9").Value = "'0.3261"
$ws.Range("E9").Value = "'  -0.33%  "
$ws.Range("D10").Value = "'1.148"
$ws.Range("D11").Value = "'0.07098"
$ws.Range("E11").Value = "'  +2.45%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "'  -0.01%  "
$ws.Range("D13").Value = "'6.105"
$ws.Range("E13").Value = "'  +3.05%  "
$ws.Range("D14").Value = "'19.70"
$ws.Range("E14").Value = "'  +2.52%  "
$ws.Range("D15").Value = "'1.674.98"
$ws.Range("E15").Value = "'  +1.68%  "
$ws.Range("D16").Value = "'6.658"
$ws.Range("E16").Value = "'  +1.07%  "
$ws.Range("D17").Value = "'0.00001052"
$ws.Range("E17").Value = "'  +1.45%  "
$ws.Range("D18").Value = "'0.06588"
$ws.Range("E18").Value = "'  +1.30%  "
$ws.Range("D19").Value = "'0.9999"
$ws.Range("E19").Value = "'  +0.04%  "
$ws.Range("D20").Value = "'79.02"
$ws.Range("E20").Value = "'  +3.44%  "
$ws.Range("D21").Value = "'15.96"
$ws.Range("E21").Value = "'  +1.95%  "
$ws.Range("D22").Value = "'5.937"
$ws.Range("E22").Value = "'  +0.41%  "
$ws.Range("D23").Value = "'12.90"
$ws.Range("E23").Value = "'  +4.59%  "
$ws.Range("D24").Value = "'24.957.44"
$ws.Range("E24").Value = "'  +2.42%  "
$ws.Range("D25").Value = "'2.453"
$ws.Range("D26").Value = "'2.425"
$ws.Range("E26").Value = "'  +3.73%  "
$ws.Range("D27").Value = "'148.11"
$ws.Range("E27").Value = "'  +1.27%  "
$ws.Range("D28").Value = "'18.81"
$ws.Range("E28").Value = "'  +3.13%  "
$ws.Range("D29").Value = "'1.866.59"
$ws.Range("E29").Value = "'  +2.02%  "
$ws.Range("D30").Value = "'126.10"
$ws.Range("E30").Value = "'  +1.55%  "
$ws.Range("D31").Value = "'1.191"
$ws.Range("E31").Value = "'  +2.20%  "
$ws.Range("D32").Value = "'4.078"
$ws.Range("E32").Value = "'  +0.81%  "
$ws.Range("D33").Value = "'5.804"
$ws.Range("E33").Value = "'  +4.29%  "
$ws.Range("D34").Value = "'0.08508"
$ws.Range("E34").Value = "'  +1.76%  "
$ws.Range("E35").Value = "'  -1.60%  "
$ws.Range("D36").Value = "'12.35"
$ws.Range("E36").Value = "'  +0.99%  "
$ws.Range("D37").Value = "'5.198"
$ws.Range("E37").Value = "'  -0.16%  "
$ws.Range("D38").Value = "'0.02257"
$ws.Range("E38").Value = "'  +2.73%  "
$ws.Range("D42").Value = "'8.266"
$ws.Range("E42").Value = "'  +1.35%  "
$ws.Range("E43").Value = "'  -0.10%  "
$ws.Range("D44").Value = "'0.5963"
$ws.Range("E44").Value = "'  +2.67%  "
$ws.Range("D45").Value = "'13.69"
$ws.Range("E45").Value = "'  +9.16%  "
$ws.Range("E46").Value = "'  +3.14%  "
$ws.Range("D47").Value = "'0.5748"
$ws.Range("E47").Value = "'  +3.97%  "
$ws.Range("D48").Value = "'125.71"
$ws.Range("E48").Value = "'  +3.37%  "
$ws.Range("D49").Value = "'1.973"
$ws.Range("E49").Value = "'  +2.10%  "
$ws.Range("D50").Value = "'0.07029"
$ws.Range("E50").Value = "'  +2.05%  "
$ws.Range("D51").Value = "'1.193"
$ws.Range("E51").Value = "'  +3.43%  "

# Rows 39-41: ranking reorder (TrustWalletToken, Hedera, Algorand moved up
# past Algorand, which drops to the bottom of this block)
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'1.231"
$ws.Range("E39").Value = "'  +1.97%  "

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "'0.06039"
$ws.Range("E40").Value = "'  +0.37%  "

$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "'0.2098"
$ws.Range("E41").Value = "'  +2.79%  "
